$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# OLE colors (R + G*256 + B*65536).
$GREEN = 5296274   # RGB(146,208,80) -> FF92D050
$RED   = 255       # RGB(255,0,0)    -> FFFF0000

# ---------------------------------------------------------------------------
# 1) C5 / C6 gain a second "Granularity: ..." line and word-wrap.
#    Their fill/border/center-alignment are already correct and untouched.
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "Data Sources: Same`nGranularity: Same"
$ws.Range("C6").Value = "Data Sources: Different`nGranularity: Same"
$ws.Range("C5:C6").WrapText = $true

# ---------------------------------------------------------------------------
# 2) D5, E6 and F6 move from the old placeholder fill to the new solid
#    green fill (their border/alignment/wrap/text stay as-is).
# ---------------------------------------------------------------------------
$ws.Range("D5").Interior.Color = $GREEN
$ws.Range("E6").Interior.Color = $GREEN
$ws.Range("F6").Interior.Color = $GREEN

# ---------------------------------------------------------------------------
# 3) Rows 7 and 8 were an empty templated pair of rows; now they host a
#    second "granularity" comparison row plus an "aggregation and join" note.
# ---------------------------------------------------------------------------

# C7 / C8 reuse C5's look (theme fill, thin box border, center + wrap).
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C7:C8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("C7").Value = "Data Sources: Same `nGranularity: Different"
$ws.Range("C8").Value = "Data Sources: Different `nGranularity: Different"

# D7:E7 - solid red fill, thin box border, centered (like D6/E5).
$r = $ws.Range("D7:E7")
$r.Interior.Color = $RED
$r.Borders.LineStyle = 1
$r.Borders.Weight = 2
$r.HorizontalAlignment = -4108

# D8:E8 - solid red fill, thin box border, default (no) alignment.
$r2 = $ws.Range("D8:E8")
$r2.Interior.Color = $RED
$r2.Borders.LineStyle = 1
$r2.Borders.Weight = 2

# F7:F8 - merged cell: green fill, centered + vertically centered + wrapped,
# thin border all around the merged block but none on the internal seam.
$ws.Range("F7:F8").Merge()
$ws.Range("F7").Value = "When we have to do aggregation and join"

$full = $ws.Range("F7:F8")
$full.Interior.Color = $GREEN
$full.HorizontalAlignment = -4108
$full.VerticalAlignment = -4108
$full.WrapText = $true

$f7 = $ws.Range("F7")
$f7.Borders.LineStyle = 1
$f7.Borders.Weight = 2
$f7.Borders.Item(9).LineStyle = -4142   # no border on the seam shared with F8

$f8 = $ws.Range("F8")
$f8.Borders.LineStyle = 1
$f8.Borders.Weight = 2
$f8.Borders.Item(8).LineStyle = -4142   # no border on the seam shared with F7

# ---------------------------------------------------------------------------
# 4) Row heights grow to fit the (now) two-line / wrapped text in rows
#    5, 7 and 8 (row 6 already had its taller 31.8pt height).
# ---------------------------------------------------------------------------
$ws.Rows(5).RowHeight = 21.6
$ws.Rows(7).RowHeight = 21.6
$ws.Rows(8).RowHeight = 21.6
